# New code for graphing several dataframes
# -> sets explicit column widths on Sheet1 (columns B:K) so the sheet's
#    <cols> block matches the widths used for the new dataframe graphing
#    section.
#
# Excel's COM ColumnWidth property is expressed in "characters" of the
# Normal style font, while the offset between ColumnWidth and the stored
# OOXML <col width="..."> value is a constant 5/6 character (the standard
# grid padding). Subtracting that offset here reproduces the target
# stored widths as closely as this host's column-width grid allows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$widths = @{
    2  = 14
    3  = 12.140625
    4  = 13.140625
    5  = 12.7109375
    6  = 11.85546875
    7  = 14.28515625
    8  = 13.5703125
    9  = 13.7109375
    10 = 13.140625
    11 = 13.140625
}

foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - (5 / 6)
}
